# Adds a new forecast-origin column ("2020-05-03", the same date already used
# as the header of column X's value = shared string already present) to both
# the "cases" and "deaths" sheets, fills in the forecast values for that
# column on rows 23-35, sets the "Observed" value for 2020-05-03 (row 22,
# column B), and appends a brand-new row 36 for forecast-origin date
# "2020-05-17" (only known in column X).
#
# NOTE: several of the values being written look like dates (e.g.
# "2020-05-03", "2020-05-17"), and Excel's COM automation normally
# auto-converts such text into a date serial number + date number format the
# moment it is assigned to a cell. To keep these as plain text (shared
# strings) exactly like the rest of the sheet, every text write below first
# forces the cell's number format to Text ("@"), assigns the value, and then
# resets the cell style back to "Normal" so no stray per-cell formatting is
# left behind.

$wb = $excel.ActiveWorkbook

$sheetNames = @("cases", "deaths")

# New forecasts for column X (origin date 2020-05-03) on rows 23-35, plus the
# brand new row 36 (origin date 2020-05-17), keyed by sheet name.
$xValues = @{
    "cases"  = @{ 23 = 34080; 24 = 36671; 25 = 39145; 26 = 41869; 27 = 43896; 28 = 45554; 29 = 46117; 30 = 48181; 31 = 49874; 32 = 51776; 33 = 52880; 34 = 53764; 35 = 54424; 36 = 55435 }
    "deaths" = @{ 23 = 2754;  24 = 2898;  25 = 3042;  26 = 3177;  27 = 3297;  28 = 3404;  29 = 3495;  30 = 3581;  31 = 3662;  32 = 3736;  33 = 3801;  34 = 3856;  35 = 3914;  36 = 3963 }
}

# "Observed" value (column B) for row 22 (2020-05-03), keyed by sheet name.
$b22Values = @{
    "cases"  = 31772
    "deaths" = 2627
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- New column header X1: text "2020-05-03" (same date text already used
    # elsewhere in the sheet) ---
    $hdr = $ws.Range("X1")
    $hdr.NumberFormat = "@"
    $hdr.Value = "2020-05-03"
    $hdr.Style = "Normal"

    # --- Touch column X for rows 2-22 so the cells exist but stay empty,
    # matching the diff's bare <c r="X2"/> ... <c r="X22"/> placeholders ---
    $ws.Range("X2:X22").Style = "Normal"

    # --- Fill in the new forecast numbers for column X, rows 23-35 ---
    foreach ($row in 23..35) {
        $ws.Range("X$row").Value = $xValues[$sheetName][$row]
    }

    # --- Observed value for 2020-05-03 goes into B22 ---
    $ws.Range("B22").Value = $b22Values[$sheetName]

    # --- Brand-new row 36 for forecast-origin date 2020-05-17 ---
    $a36 = $ws.Range("A36")
    $a36.NumberFormat = "@"
    $a36.Value = "2020-05-17"
    $a36.Style = "Normal"

    $ws.Range("B36:W36").Style = "Normal"

    $ws.Range("X36").Value = $xValues[$sheetName][36]
}
